$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New rows: 10, 11, 12 -> columns A (번호), B (문제 이름), C (틀린 날짜)
$rows = @(
    @{ Row = 10; No = 17071; Name = "숨바꼭질 5";   Date = "2025-09-04" },
    @{ Row = 11; No = 16637; Name = "괄호추가하기"; Date = "2025-09-04" },
    @{ Row = 12; No = 3197;  Name = "백조의 호수";  Date = "2025-09-05" }
)

$ws.Range("C9").Copy()
$ws.Range("C10:C12").PasteSpecial(-4122)

foreach ($r in $rows) {
    $ws.Cells.Item($r.Row, 1).Value = $r.No
    $ws.Cells.Item($r.Row, 2).Value = $r.Name
    $ws.Cells.Item($r.Row, 3).Value = $r.Date
}

$ws.Range("D12").Select()
